$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# The "形態" (data type) for the CreateDate and LastUpdate fields changes
# from DATE to TIMESTAMP.
$ws.Range("D12").Value = "TIMESTAMP"
$ws.Range("D14").Value = "TIMESTAMP"

# Reflect the author's last selection on this sheet.
$ws.Range("D14").Select()
